# Applies the cryptos.xlsx cell-value updates described in the commit diff
# (cryptocurrency price/volume refresh + two row re-orderings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.422.64"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "2.947.43"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'570.95"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").Value = "'161.24"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").Value = "2.942.79"
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("D10").Value = "'6.67"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("E11").Value = "  -4.43%  "
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").Value = "'0.0000243"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").Value = "'34.73"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "65.538.57"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "3.440.04"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "'7.06"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").Value = "2.947.41"
$ws.Range("D20").Value = "'15.81"
$ws.Range("E20").Value = "  +13.34%  "
$ws.Range("D21").Value = "'445.63"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("D22").Value = "'0.694"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").Value = "'7.28"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").Value = "'82.05"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("D26").Value = "'12.24"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.01"
$ws.Range("E28").Value = "  -5.83%  "
$ws.Range("D29").Value = "'2.50"
$ws.Range("E29").Value = "  +6.81%  "
$ws.Range("D30").Value = "'8.03"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("D31").Value = "'2.59"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "'0.0000101"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("E33").Value = "  +3.55%  "
$ws.Range("D34").Value = "'27.13"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'0.970"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("E37").Value = "  -1.65%  "
$ws.Range("D38").Value = "'45.98"
$ws.Range("E38").Value = "  +5.30%  "
$ws.Range("D39").Value = "'49.03"
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").Value = "'1.97"
$ws.Range("E40").Value = "  -7.61%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.302"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.121"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  -6.35%  "
$ws.Range("D44").Value = "'8.50"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "'381.59"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("D47").Value = "2.678.72"
$ws.Range("E47").Value = "  -4.20%  "
$ws.Range("D48").Value = "'133.16"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "'23.78"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "'2.16"
$ws.Range("E51").Value = "  +1.56%  "
